# "solved sort the jumbled numbers"
# Append a new tracker row (row 24) for the LeetCode problem
# "Sort the Jumbled Numbers", mirroring the existing rows' layout:
#   A=number, B=difficulty, C=Question, D=url, E=Approach,
#   F=efficient O, G=highlight, H=date

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24

$ws.Range("A$newRow").Value = 2191
$ws.Range("B$newRow").Value = "Medium"
$ws.Range("C$newRow").Value = "Sort the Jumbled Numbers"
$ws.Range("D$newRow").Value = "https://leetcode.com/problems/sort-the-jumbled-numbers/description/"
$ws.Range("E$newRow").Value = "Array"
$ws.Range("F$newRow").Value = "O(nd+logn)"
$ws.Range("G$newRow").Value = "Encode the numbers according to the map, store (encoded_num, index) in a list. Exploit python's sort feature."
$ws.Range("H$newRow").Value = 45497

# Turn the url cell into a real hyperlink, same as every other row.
$ws.Hyperlinks.Add($ws.Range("D$newRow"), "https://leetcode.com/problems/sort-the-jumbled-numbers/description/")

# Match the look of the row above (fill colour, borders, hyperlink font,
# date number format, etc.) by copying its formatting down.
$ws.Range("A23:H23").Copy()
$ws.Range("A24:H24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reflect where the user was working after adding the new entry.
$ws.Range("D29").Select() | Out-Null
